$d = $word.ActiveDocument

# Update the date heading
$d.Content.Find.Execute("2023-09-17 Sunday", $true, $false, $false, $false, $false, $true, 1, $false, "2023-09-18 Monday", 2) | Out-Null

$t = $d.Tables.Item(1)

$t.Cell(1,1).Range.Text = "89÷2=44, 1"
$t.Cell(1,2).Range.Text = "98÷3=32, 2"
$t.Cell(1,3).Range.Text = "82÷9=9, 1"
$t.Cell(1,4).Range.Text = "52÷2=26, 0"
$t.Cell(1,5).Range.Text = "91÷3=30, 1"

$t.Cell(5,1).Range.Text = "92÷4=23, 0"
$t.Cell(5,2).Range.Text = "46÷2=23, 0"
$t.Cell(5,3).Range.Text = "92÷7=13, 1"
$t.Cell(5,4).Range.Text = "41÷4=10, 1"
$t.Cell(5,5).Range.Text = "65÷2=32, 1"

$t.Cell(9,1).Range.Text = "55÷5=11, 0"
$t.Cell(9,2).Range.Text = "73÷7=10, 3"
$t.Cell(9,3).Range.Text = "53÷7=7, 4"
$t.Cell(9,4).Range.Text = "47÷7=6, 5"
$t.Cell(9,5).Range.Text = "59÷2=29, 1"

$t.Cell(13,1).Range.Text = "22÷8=2, 6"
$t.Cell(13,2).Range.Text = "34÷4=8, 2"
$t.Cell(13,3).Range.Text = "53÷9=5, 8"
$t.Cell(13,4).Range.Text = "37÷6=6, 1"
$t.Cell(13,5).Range.Text = "57÷8=7, 1"

$t.Cell(17,1).Range.Text = "68÷4=17, 0"
$t.Cell(17,2).Range.Text = "39÷7=5, 4"
$t.Cell(17,3).Range.Text = "80÷3=26, 2"
$t.Cell(17,4).Range.Text = "48÷4=12, 0"
$t.Cell(17,5).Range.Text = "98÷9=10, 8"
